$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values must stay as literal text (not be reinterpreted as
# numbers/dates by Excel), so force Text number format before assigning.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.241.02'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.906.50'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.51'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5377'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3821'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07300'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.27'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9061'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08201'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '95.70'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.361'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008657'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '27.266.89'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.048'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.055.94'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.78'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.519'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '149.06'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.306'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.748'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '116.80'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.830'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.723'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09225'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8288'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05091'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.003'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.321'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.678'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5935'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.078'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.404'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.665'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '116.82'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5096'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.17'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.646'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06158'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.63'

# Coin name / link / volume columns are safe to assign directly as text.
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("E7").Value = '  +2.94%  '
$ws.Range("E8").Value = '  +1.83%  '
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("E10").Value = '  +5.19%  '
$ws.Range("E11").Value = '  +0.81%  '
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("E13").Value = '  -1.17%  '
$ws.Range("E14").Value = '  +1.65%  '
$ws.Range("E15").Value = '  -0.20%  '
$ws.Range("E16").Value = '  +2.40%  '
$ws.Range("E17").Value = '  +0.65%  '
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("E20").Value = '  -0.66%  '
$ws.Range("E21").Value = '  -44.31%  '
$ws.Range("E22").Value = '  +0.88%  '
$ws.Range("E23").Value = '  +1.89%  '
$ws.Range("E24").Value = '  +0.93%  '
$ws.Range("E25").Value = '  +0.75%  '
$ws.Range("E26").Value = '  +1.17%  '
$ws.Range("E27").Value = '  +0.17%  '
$ws.Range("E28").Value = '  +1.54%  '
$ws.Range("E29").Value = '  +0.56%  '
$ws.Range("E30").Value = '  -3.84%  '
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("E32").Value = '  +4.46%  '
$ws.Range("E33").Value = '  +0.98%  '
$ws.Range("E34").Value = '  +0.26%  '
$ws.Range("E35").Value = '  +1.80%  '
$ws.Range("E36").Value = '  -3.42%  '
$ws.Range("E37").Value = '  +4.49%  '
$ws.Range("E38").Value = '  +5.02%  '
$ws.Range("E39").Value = '  +0.96%  '
$ws.Range("E40").Value = '  +0.36%  '
$ws.Range("E41").Value = '  +5.29%  '
$ws.Range("E42").Value = '  +2.09%  '
$ws.Range("E43").Value = '  +1.45%  '
$ws.Range("E44").Value = '  +4.52%  '
$ws.Range("E45").Value = '  +0.88%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("E47").Value = '  +0.70%  '
$ws.Range("E48").Value = '  +1.85%  '
$ws.Range("E49").Value = '  +0.81%  '
